$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.723.67"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.66"
$ws.Range("E3").Value = "  -1.19%  "

$ws.Range("E4").Value = "  -2.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.63"
$ws.Range("E5").Value = "  -1.77%  "

$ws.Range("E6").Value = "  -1.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4315"
$ws.Range("E7").Value = "  -2.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3753"
$ws.Range("E8").Value = "  -2.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07360"
$ws.Range("E9").Value = "  -1.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8785"
$ws.Range("E10").Value = "  -1.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.67"
$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.858.91"
$ws.Range("E12").Value = "  -0.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.739"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.451"
$ws.Range("E14").Value = "  -2.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07150"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.04"
$ws.Range("E16").Value = "  +4.49%  "

$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009009"
$ws.Range("E18").Value = "  -1.51%  "

$ws.Range("E19").Value = "  -2.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.53"
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.736.38"
$ws.Range("E21").Value = "  -0.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.223"
$ws.Range("E22").Value = "  -2.31%  "

$ws.Range("E23").Value = "  -2.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.080.46"
$ws.Range("E24").Value = "  -1.17%  "

$ws.Range("E25").Value = "  -1.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.51"
$ws.Range("E26").Value = "  -2.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.67"
$ws.Range("E27").Value = "  -1.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.184"
$ws.Range("E28").Value = "  +9.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.385"
$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.14"
$ws.Range("E30").Value = "  +0.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08947"
$ws.Range("E31").Value = "  -1.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.234"
$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7798"
$ws.Range("E33").Value = "  -0.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.570"
$ws.Range("E34").Value = "  -1.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.919"
$ws.Range("E35").Value = "  -3.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.013"
$ws.Range("E36").Value = "  -2.01%  "

$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05371"
$ws.Range("E38").Value = "  -0.25%  "

$ws.Range("E39").Value = "  -0.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.323"
$ws.Range("E40").Value = "  +5.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.900"
$ws.Range("E41").Value = "  +0.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5145"
$ws.Range("E42").Value = "  -1.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1692"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.840"
$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.78"
$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "108.72"
$ws.Range("E46").Value = "  -3.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4803"
$ws.Range("E47").Value = "  +0.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06477"
$ws.Range("E48").Value = "  -2.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.697"
$ws.Range("E49").Value = "  -2.04%  "

$ws.Range("E50").Value = "  -2.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.856"
$ws.Range("E51").Value = "  -3.59%  "
